# Apply updated dSF (column F) values for specific rows, as part of the
# "repull data, push all data, mean calculation" correction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    16 = 3
    31 = 1
    33 = 1
    37 = 2
    40 = 0
    41 = -1
    42 = 0
    44 = -1
    53 = 2
    58 = -1
    63 = 2
    67 = 0
    70 = -2
    76 = 0
    81 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
